$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear leftover formatting on C2:D2 (reset to default style)
$ws.Range("C2:D2").ClearFormats()

# Update Programming Status for Week2 and Week3 to "In Progress"
$ws.Range("D4").Value = "In Progress"
$ws.Range("D5").Value = "In Progress"

# Update Slide Status for Week10 (Collections) to "Done" - also match the
# green fill used by the other "Done" cells (e.g. C3) instead of the
# yellow "In Progress" fill it previously had
$ws.Range("C12").Value = "Done"
$ws.Range("C12").Interior.Color = $ws.Range("C3").Interior.Color

# Move active cell selection to F5
$ws.Range("F5").Select()
